# Add the "Other Perturbations" section to the manifest sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New section header (bold, like the other section headers on this sheet)
$ws.Range("A24").Value = "[Other Perturbations]"
$ws.Range("A24").Font.Bold = $true

# New field row under the section header
$ws.Range("A25").Value = "#Description"

# Leave the active cell/selection on the newly added row, as in the target workbook
$null = $ws.Range("A25").Select()
